# This script swaps the contents of rows 10 and 11 (data rows for
# "CONNECT GLOBAL BUSINESS COMPANY Ltd" and "LU SUPPLY COMPANY  LTD")
# on the active worksheet, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually contain data in rows 10 and 11 (L and M are blank).
# Column R ("bids won") is a genuine number and must stay numeric; all the
# others hold text in the source data (ids/phone numbers/account numbers
# included), even when their contents look like plain digit strings.
$textColumns = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "N", "O", "P", "Q")
$numericColumns = @("R")

foreach ($col in $textColumns) {
    $cell10 = $ws.Range($col + "10")
    $cell11 = $ws.Range($col + "11")

    # NOTE: reading the `.Value` property as an expression on this runtime
    # returns a descriptor string instead of the cell's contents (writes via
    # `.Value = ...` work correctly though). `.Value2` also silently
    # coerces all-digit text into a Double (e.g. for the long national-id
    # numbers), which loses the exact literal digits once you need to
    # re-stringify it. `.Text` gives back the literal characters that were
    # stored, so use that for reads here.
    $value10 = [string]$cell10.Text
    $value11 = [string]$cell11.Text

    # The source data keeps every one of these fields as plain text (even
    # the all-digit ids/phone numbers/bank account numbers, some of which
    # have significant leading zeros). Writing a bare digit-string back
    # through `.Value` lets Excel reinterpret it as a number, dropping
    # leading zeros, so force those via a leading apostrophe to keep them
    # as text - exactly like typing `'0788771482` into a cell.
    if ($value11 -match '^[0-9]+$') {
        $cell10.Value = "'" + $value11
    } else {
        $cell10.Value = $value11
    }

    if ($value10 -match '^[0-9]+$') {
        $cell11.Value = "'" + $value10
    } else {
        $cell11.Value = $value10
    }
}

foreach ($col in $numericColumns) {
    $cell10 = $ws.Range($col + "10")
    $cell11 = $ws.Range($col + "11")

    $value10 = $cell10.Value2
    $value11 = $cell11.Value2

    $cell10.Value = $value11
    $cell11.Value = $value10
}
